$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2338.44
$ws.Range("I113").Value = 2212.6191
$ws.Range("J113").Value = 2999
$ws.Range("K113").Value = 2212.6191
$ws.Range("L113").Value = 2999
$ws.Range("M113").Value = 1041.3809
$ws.Range("N113").Value = -9507
$ws.Range("H133").Value = 16896
$ws.Range("J133").Value = 16896
$ws.Range("L133").Value = 16896
$ws.Range("N133").Value = -27016

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1668.8334
$ws.Range("I2").Value = 1450
$ws.Range("J2").Value = 2106.5
$ws.Range("K2").Value = 1450
$ws.Range("L2").Value = 2106.5
$ws.Range("M2").Value = -1337
$ws.Range("N2").Value = -2332.5
$ws.Range("H45").Value = 5684644.5
$ws.Range("I45").Value = 7578626.5
$ws.Range("J45").Value = 2699
$ws.Range("K45").Value = 7578626.5
$ws.Range("L45").Value = 2699
$ws.Range("M45").Value = -7578249.5
$ws.Range("N45").Value = -3453
$ws.Range("H61").Value = 2063.7708
$ws.Range("I61").Value = 2128.4634
$ws.Range("J61").Value = 1684.8572
$ws.Range("K61").Value = 2128.4634
$ws.Range("L61").Value = 1684.8572
$ws.Range("M61").Value = -1916.4634
$ws.Range("N61").Value = -2108.8572
$ws.Range("H74").Value = 1582.5161
$ws.Range("I74").Value = 1274.1904
$ws.Range("J74").Value = 2230
$ws.Range("K74").Value = 1274.1904
$ws.Range("L74").Value = 2230
$ws.Range("M74").Value = -400.1904
$ws.Range("N74").Value = -3978
$ws.Range("H77").Value = 1582.5161
$ws.Range("I77").Value = 1274.1904
$ws.Range("J77").Value = 2230
$ws.Range("K77").Value = 6370.951999999999
$ws.Range("L77").Value = 11150
$ws.Range("M77").Value = -2002.951999999999
$ws.Range("N77").Value = -19886
$ws.Range("H102").Value = 1212.5
$ws.Range("I102").Value = 1100
$ws.Range("J102").Value = 1400
$ws.Range("K102").Value = 1100
$ws.Range("L102").Value = 1400
$ws.Range("M102").Value = 522
$ws.Range("N102").Value = -4644
$ws.Range("H107").Value = 28818.666
$ws.Range("J107").Value = 28818.666
$ws.Range("L107").Value = 28818.666
$ws.Range("N107").Value = -36498.666
$ws.Range("H109").Value = 32415.777
$ws.Range("J109").Value = 32415.777
$ws.Range("L109").Value = 32415.777
$ws.Range("N109").Value = -35189.777
$ws.Range("H116").Value = 1668.8334
$ws.Range("I116").Value = 1450
$ws.Range("J116").Value = 2106.5
$ws.Range("K116").Value = 1450
$ws.Range("L116").Value = 2106.5
$ws.Range("M116").Value = 844
$ws.Range("N116").Value = -6694.5
$ws.Range("H132").Value = 6758884.5
$ws.Range("I132").Value = 10418440
$ws.Range("J132").Value = 2781.6924
$ws.Range("K132").Value = 31255320
$ws.Range("L132").Value = 8345.0772
$ws.Range("M132").Value = -31252790
$ws.Range("N132").Value = -13405.0772
$ws.Range("H136").Value = 2063.7708
$ws.Range("I136").Value = 2128.4634
$ws.Range("J136").Value = 1684.8572
$ws.Range("K136").Value = 6385.3902
$ws.Range("L136").Value = 5054.571599999999
$ws.Range("M136").Value = -3835.3902
$ws.Range("N136").Value = -10154.5716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1668.8334
$ws.Range("I3").Value = 1450
$ws.Range("J3").Value = 2106.5
$ws.Range("K3").Value = 1450
$ws.Range("L3").Value = 2106.5
$ws.Range("M3").Value = -1336
$ws.Range("N3").Value = -2334.5
$ws.Range("H20").Value = 1802.6666
$ws.Range("I20").Value = 1903.5927
$ws.Range("J20").Value = 1575.5834
$ws.Range("K20").Value = 1903.5927
$ws.Range("L20").Value = 1575.5834
$ws.Range("M20").Value = -1656.5927
$ws.Range("N20").Value = -2069.5834
$ws.Range("H105").Value = 3396.721
$ws.Range("I105").Value = 1458.6
$ws.Range("J105").Value = 4435
$ws.Range("K105").Value = 1458.6
$ws.Range("L105").Value = 4435
$ws.Range("M105").Value = 288.4000000000001
$ws.Range("N105").Value = -7929
$ws.Range("H134").Value = 4018.4583
$ws.Range("I134").Value = 2845.25
$ws.Range("K134").Value = 8535.75
$ws.Range("M134").Value = -6000.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2651.4614
$ws.Range("I31").Value = 1917.8077
$ws.Range("J31").Value = 3385.1155
$ws.Range("K31").Value = 1917.8077
$ws.Range("L31").Value = 3385.1155
$ws.Range("M31").Value = -1622.8077
$ws.Range("N31").Value = -3975.1155
$ws.Range("H34").Value = 2651.4614
$ws.Range("I34").Value = 1917.8077
$ws.Range("J34").Value = 3385.1155
$ws.Range("K34").Value = 1917.8077
$ws.Range("L34").Value = 3385.1155
$ws.Range("M34").Value = -1715.8077
$ws.Range("N34").Value = -3789.1155
$ws.Range("H58").Value = 1796.3513
$ws.Range("I58").Value = 857.4211
$ws.Range("J58").Value = 2787.4443
$ws.Range("K58").Value = 857.4211
$ws.Range("L58").Value = 2787.4443
$ws.Range("M58").Value = -654.4211
$ws.Range("N58").Value = -3193.4443
$ws.Range("H136").Value = 1796.3513
$ws.Range("I136").Value = 857.4211
$ws.Range("J136").Value = 2787.4443
$ws.Range("K136").Value = 2572.2633
$ws.Range("L136").Value = 8362.332900000001
$ws.Range("M136").Value = -22.26330000000007
$ws.Range("N136").Value = -13462.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3349.875
$ws.Range("I80").Value = 1899
$ws.Range("J80").Value = 3557.1428
$ws.Range("K80").Value = 5697
$ws.Range("L80").Value = 10671.4284
$ws.Range("M80").Value = -4761
$ws.Range("N80").Value = -12543.4284
$ws.Range("H83").Value = 3349.875
$ws.Range("I83").Value = 1899
$ws.Range("J83").Value = 3557.1428
$ws.Range("K83").Value = 17091
$ws.Range("L83").Value = 32014.2852
$ws.Range("M83").Value = -12411
$ws.Range("N83").Value = -41374.2852
$ws.Range("H92").Value = 371.42856
$ws.Range("I92").Value = 247.5
$ws.Range("J92").Value = 421
$ws.Range("K92").Value = 742.5
$ws.Range("L92").Value = 1263
$ws.Range("M92").Value = 505.5
$ws.Range("N92").Value = -3759

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2864.0715
$ws.Range("I102").Value = 3647
$ws.Range("J102").Value = 1654.091
$ws.Range("K102").Value = 3647
$ws.Range("L102").Value = 1654.091
$ws.Range("M102").Value = -2025
$ws.Range("N102").Value = -4898.091
$ws.Range("H126").Value = 2843.8928
$ws.Range("I126").Value = 1991.909
$ws.Range("J126").Value = 3395.1765
$ws.Range("K126").Value = 5975.727000000001
$ws.Range("L126").Value = 10185.5295
$ws.Range("M126").Value = -3505.727000000001
$ws.Range("N126").Value = -15125.5295

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10840.471
$ws.Range("I7").Value = 9831.666999999999
$ws.Range("J7").Value = 11390.728
$ws.Range("K7").Value = 9831.666999999999
$ws.Range("L7").Value = 11390.728
$ws.Range("M7").Value = -9719.666999999999
$ws.Range("N7").Value = -11614.728
$ws.Range("H16").Value = 1488.5483
$ws.Range("J16").Value = 1850
$ws.Range("L16").Value = 1850
$ws.Range("N16").Value = -2190
$ws.Range("H55").Value = 261.8095
$ws.Range("I55").Value = 136.4
$ws.Range("J55").Value = 301
$ws.Range("K55").Value = 136.4
$ws.Range("L55").Value = 301
$ws.Range("M55").Value = 36.59999999999999
$ws.Range("N55").Value = -647
$ws.Range("H122").Value = 6286.9287
$ws.Range("I122").Value = 4924.1113
$ws.Range("J122").Value = 8740
$ws.Range("K122").Value = 14772.3339
$ws.Range("L122").Value = 26220
$ws.Range("M122").Value = -12322.3339
$ws.Range("N122").Value = -31120
$ws.Range("H126").Value = 10840.471
$ws.Range("I126").Value = 9831.666999999999
$ws.Range("J126").Value = 11390.728
$ws.Range("K126").Value = 29495.001
$ws.Range("L126").Value = 34172.18399999999
$ws.Range("M126").Value = -27025.001
$ws.Range("N126").Value = -39112.18399999999
$ws.Range("H136").Value = 4947.718
$ws.Range("I136").Value = 2141.3142
$ws.Range("K136").Value = 6423.942599999999
$ws.Range("M136").Value = -3873.942599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1459.6222
$ws.Range("I136").Value = 918.2727
$ws.Range("J136").Value = 2948.3333
$ws.Range("K136").Value = 2754.8181
$ws.Range("L136").Value = 8844.999899999999
$ws.Range("M136").Value = -204.8181
$ws.Range("N136").Value = -13944.9999
